$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shrink the used range from 28 rows to 26 rows by deleting the two trailing rows
$ws.Rows(28).Delete()
$ws.Rows(27).Delete()

# Row 8
$ws.Range("A8").Value = ''
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 28
$ws.Range("D8").Value = '1.0'
$ws.Range("E8").Value = 'Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it''s ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = '0.00'
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = ''

# Row 9
$ws.Range("A9").Value = 'P. point'
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 37
$ws.Range("D9").Value = '6'
$ws.Range("E9").Value = 'On board'
$ws.Range("F9").Value = 136
$ws.Range("G9").Value = '5032.00'
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = ''

# Row 10
$ws.Range("A10").Value = 'Each'
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 16
$ws.Range("D10").Value = '9.0'
$ws.Range("E10").Value = 'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F10").Value = 219
$ws.Range("G10").Value = '3504.00'
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = ''

# Row 11
$ws.Range("A11").Value = 'Each'
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 34
$ws.Range("D11").Value = '10.0'
$ws.Range("E11").Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F11").Value = 303
$ws.Range("G11").Value = '10302.00'
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = ''

# Row 12
$ws.Range("A12").Value = ''
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 65
$ws.Range("D12").Value = '11.0'
$ws.Range("E12").Value = 'S&F following sizes (dia.) of ISI marked virgin material MMS ( IS:9537 P - III ) PVC conduit along with  ISI marked (IS:3419-1988) accessories as required  in  recess  including  cutting the wall, covering conduit and making good the same as required. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = '0.00'
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = ''

# Row 13
$ws.Range("A13").Value = 'R. mtr.'
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 31
$ws.Range("D13").Value = '17'
$ws.Range("E13").Value = '25 mm'
$ws.Range("F13").Value = 56
$ws.Range("G13").Value = '1736.00'
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = ''

# Row 14
$ws.Range("A14").Value = 'Mtr.'
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 78
$ws.Range("D14").Value = '19'
$ws.Range("E14").Value = '2 x 2.5 sq. mm. + 1x1.5sqmm'
$ws.Range("F14").Value = 81
$ws.Range("G14").Value = '6318.00'
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = ''

# Row 15
$ws.Range("A15").Value = 'Mtr.'
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 37
$ws.Range("D15").Value = '20'
$ws.Range("E15").Value = '2 x 4.0 sq. mm. + 1 x 2.5 sq. mm.'
$ws.Range("F15").Value = 122
$ws.Range("G15").Value = '4514.00'
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = ''

# Row 16
$ws.Range("A16").Value = 'Set'
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 75
$ws.Range("D16").Value = '13.0'
$ws.Range("E16").Value = 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
$ws.Range("F16").Value = 5733
$ws.Range("G16").Value = '429975.00'
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = ''

# Row 17
$ws.Range("A17").Value = 'Mtr.'
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 90
$ws.Range("D17").Value = '23'
$ws.Range("E17").Value = '8 SWG G.I. ( Hot Dipped  ) Wire '
$ws.Range("F17").Value = 20
$ws.Range("G17").Value = '1800.00'
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = ''

# Row 18
$ws.Range("A18").Value = ''
$ws.Range("B18").Value = 0
$ws.Range("C18").Value = 39
$ws.Range("D18").Value = '16.0'
$ws.Range("E18").Value = 'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = '0.00'
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = ''

# Row 19
$ws.Range("A19").Value = ''
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = '17.0'
$ws.Range("E19").Value = 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = '0.00'
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = ''

# Row 20
$ws.Range("A20").Value = 'Each'
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 100
$ws.Range("D20").Value = '35'
$ws.Range("E20").Value = '8 Way (8+2)'
$ws.Range("F20").Value = 2184
$ws.Range("G20").Value = '218400.00'
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = ''

# Row 21
$ws.Range("A21").Value = '%'
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 78
$ws.Range("D21").Value = '37'
$ws.Range("E21").Value = 'Add Tender Premium '
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = '0.00'
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = ''

# Row 22
$ws.Range("A22").Value = ''
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 81
$ws.Range("D22").Value = '38'
$ws.Range("E22").Value = 'Grand Total'
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = '0.00'
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = ''

# Row 23
$ws.Range("B23:I23").ClearContents()
$ws.Range("A23").Value = ''

# Row 24
$ws.Range("A24").Value = ''
$ws.Range("B24").Value = ''
$ws.Range("C24").Value = ''
$ws.Range("D24").Value = ''
$ws.Range("E24").Value = 'Grand Total Rs.'
$ws.Range("F24").Value = ''
$ws.Range("G24").Value = '681581.00'
$ws.Range("H24").Value = '681581.00'
$ws.Range("I24").Value = ''

# Row 25
$ws.Range("A25").Value = ''
$ws.Range("B25").Value = ''
$ws.Range("C25").Value = ''
$ws.Range("D25").Value = ''
$ws.Range("E25").Value = 'Tender Premium @ 0%'
$ws.Range("F25").Value = ''
$ws.Range("G25").Value = '0.00'
$ws.Range("H25").Value = '0.00'
$ws.Range("I25").Value = ''

# Row 26
$ws.Range("A26").Value = ''
$ws.Range("B26").Value = ''
$ws.Range("C26").Value = ''
$ws.Range("D26").Value = ''
$ws.Range("E26").Value = 'NET PAYABLE AMOUNT Rs.'
$ws.Range("F26").Value = ''
$ws.Range("G26").Value = '681581.00'
$ws.Range("H26").Value = '681581.00'
$ws.Range("I26").Value = ''

